# Apply the Mon Aug 28 19:29:48 UTC 2023 cryptos-list refresh (GitHub Actions job).
# Prices/volumes are scraped text, not numbers, so every write below targets the
# sheet as literal text -- "Price" cells such as "1.005" or "0.06330" would
# otherwise be auto-coerced by Excel into doubles and lose their trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the plain-decimal Price cells as Text so Excel keeps them as strings.
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D15",
    "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27",
    "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38",
    "D39", "D40", "D42", "D43", "D46", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.103.48"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").Value = "1.650.12"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "218.77"
$ws.Range("E5").Value = "  -0.83%  "

$ws.Range("D6").Value = "0.5235"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "0.2657"
$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("D9").Value = "0.06330"
$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("D10").Value = "20.56"
$ws.Range("E10").Value = "  -1.66%  "

$ws.Range("D11").Value = "0.07716"
$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("D12").Value = "4.582"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "1.690.30"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").Value = "1.882.57"
$ws.Range("E14").Value = "  -0.79%  "

$ws.Range("D15").Value = "0.5603"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "0.0₅8166"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "65.22"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").Value = "26.129.35"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("D20").Value = "4.684"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "191.23"
$ws.Range("E21").Value = "  -3.76%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "10.35"
$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("D23").Value = "5.970"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "144.34"
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("D26").Value = "0.1203"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("D27").Value = "7.234"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").Value = "15.90"
$ws.Range("E28").Value = "  -1.89%  "

$ws.Range("D29").Value = "1.488"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").Value = "0.05604"
$ws.Range("E30").Value = "  -5.10%  "

$ws.Range("D31").Value = "1.274"

$ws.Range("D32").Value = "3.488"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("D33").Value = "3.362"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").Value = "1.573"
$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "2.792"
$ws.Range("E35").Value = "  -1.52%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.9473"
$ws.Range("E36").Value = "  -1.60%  "

$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").Value = "0.5729"
$ws.Range("E38").Value = "  -1.62%  "

$ws.Range("D39").Value = "0.01590"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").Value = "5.983"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D42").Value = "0.8404"
$ws.Range("E42").Value = "  -2.10%  "

$ws.Range("D43").Value = "101.46"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "1.014.17"
$ws.Range("E44").Value = "  -5.88%  "

$ws.Range("D45").Value = "1.793.56"
$ws.Range("E45").Value = "  -0.78%  "

$ws.Range("D46").Value = "58.22"
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").Value = "0.05330"
$ws.Range("E48").Value = "  +3.46%  "

$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.4347"
$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.013"
$ws.Range("E51").Value = "  -0.76%  "

